$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 641 (shifts old 641..724 down to 643..726)
$ws.Rows.Item(641).Resize(2).Insert()

# New row 641: Primera
$ws.Cells.Item(641, 1).Value = 3
$ws.Cells.Item(641, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(641, 3).Value = "Coquimbo"
$ws.Cells.Item(641, 4).Value = 44776
$ws.Cells.Item(641, 5).Value = 5
$ws.Cells.Item(641, 6).Value = 100112008
$ws.Cells.Item(641, 7).Value = "Coliflor"
$ws.Cells.Item(641, 8).Value = "Sin especificar"
$ws.Cells.Item(641, 9).Value = "Primera"
$ws.Cells.Item(641, 10).Value = 2150
$ws.Cells.Item(641, 11).Value = 1000
$ws.Cells.Item(641, 12).Value = 1100
$ws.Cells.Item(641, 13).Value = 1044
$ws.Cells.Item(641, 14).Value = "`$/unidad"
$ws.Cells.Item(641, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(641, 16).Value = 1044
$ws.Cells.Item(641, 17).Value = 1
$ws.Cells.Item(641, 18).Value = "Hortaliza"

# New row 642: Segunda
$ws.Cells.Item(642, 1).Value = 3
$ws.Cells.Item(642, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(642, 3).Value = "Coquimbo"
$ws.Cells.Item(642, 4).Value = 44776
$ws.Cells.Item(642, 5).Value = 5
$ws.Cells.Item(642, 6).Value = 100112008
$ws.Cells.Item(642, 7).Value = "Coliflor"
$ws.Cells.Item(642, 8).Value = "Sin especificar"
$ws.Cells.Item(642, 9).Value = "Segunda"
$ws.Cells.Item(642, 10).Value = 900
$ws.Cells.Item(642, 11).Value = 800
$ws.Cells.Item(642, 12).Value = 800
$ws.Cells.Item(642, 13).Value = 800
$ws.Cells.Item(642, 14).Value = "`$/unidad"
$ws.Cells.Item(642, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(642, 16).Value = 800
$ws.Cells.Item(642, 17).Value = 1
$ws.Cells.Item(642, 18).Value = "Hortaliza"
